# Update the "想去人数" (want-to-go count) numbers on the "展览" and "全部类型"
# sheets to reflect the newly generated data (F2: 124 -> 127, F3: 29 -> 31).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 127
    $ws.Range("F3").Value = 31
}
